$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 23:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1055455
$ws.Range("C4").Value = 19690
$ws.Range("E4").Value = 849920

# Row 8 - Reino Unido
$ws.Range("G8").Value = 765

# Row 9 - Alemania
$ws.Range("B9").Value = 161197
$ws.Range("C9").Value = 1285
$ws.Range("E9").Value = 34392
$ws.Range("G9").Value = 91
$ws.Range("H9").Value = 6405

# Row 15 - Canada
$ws.Range("B15").Value = 51248
$ws.Range("C15").Value = 1222
$ws.Range("D15").Value = 20100
$ws.Range("E15").Value = 28163
$ws.Range("G15").Value = 126
$ws.Range("H15").Value = 2985

# Row 27 - Israel
$ws.Range("B27").Value = 15834
$ws.Range("C27").Value = 106
$ws.Range("D27").Value = 8233
$ws.Range("E27").Value = 7386
$ws.Range("F27").Value = 115
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 215

# Row 147 - Bermudas
$ws.Range("B147").Value = 111
$ws.Range("C147").Value = 1
$ws.Range("D147").Value = 48
$ws.Range("E147").Value = 57
